# Assignment 2 edit script
# 1) Change the date "September 22, 2018" -> "October 8, 2018"
# 2) Change "Assignment 1" -> "Assignment 2"
# (the "_GoBack" bookmark naturally relocates to sit right after the
#  newly typed "8", matching Word's own behavior of tracking the most
#  recent edit point)

$d = $word.ActiveDocument

# --- 1) Date: "September 22, 2018 " -> "October 8, 2018 " -------------
# Locate "September 22" (the part that is actually being retyped) via Find
# so we don't depend on a hard-coded paragraph index.
$dateRange = $d.Content
$dateRange.Find.Execute("September 22")
$dateRange.Text = "October 8"

# Force the freshly-typed text to break into separate runs the same way
# Word does while you type ("October" / " " / "8"), by toggling a
# character property on/off over each growing prefix.
$r1 = $d.Range($dateRange.Start, $dateRange.Start + 7)         # "October"
$r1.Font.Bold = 1
$r1.Font.Bold = 0

$r2 = $d.Range($dateRange.Start, $dateRange.Start + 8)         # "October "
$r2.Font.Bold = 1
$r2.Font.Bold = 0

$r3 = $d.Range($dateRange.Start, $dateRange.Start + 9)         # "October 8"
$r3.Font.Bold = 1
$r3.Font.Bold = 0

# Drop the last-edit marker right after the "8", exactly where Word
# leaves it after you finish typing.
$goBack = $d.Range($dateRange.Start + 9, $dateRange.Start + 9)
$d.Bookmarks.Add("_GoBack", $goBack)

# --- 2) Assignment number -------------------------------------------
$d.Content.Find.Execute("Assignment 1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Assignment 2", 2)
